$d = $word.ActiveDocument

# Locate the manuscript title paragraph by its text content, to be robust
# against any structural drift rather than hard-coding a paragraph index.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Climatic Variables Alone*") {
        $targetPara = $p
        break
    }
}

# Fall back to the known paragraph position if the text lookup ever fails.
if ($targetPara -eq $null) {
    $targetPara = $d.Paragraphs.Item(3)
}

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14"><w:body><w:p w14:paraId="2BE3D0EE" w14:textId="549224B6" w:rsidR="00CC5925" w:rsidRPr="00CC5925" w:rsidRDefault="00CC5925" w:rsidP="00CC5925"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr></w:pPr><w:r w:rsidRPr="00CC5925"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr><w:t xml:space="preserve">Climatic Variables Alone do not Determine Ungulate Distributions in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00CC5925"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr><w:t>Afrotropics</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetPara.Range.InsertXML($xml)
